$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.209.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.17%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.835.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.24%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'241.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.04%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  -2.81%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.07375"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.63%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.2925"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.34%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -2.22%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07719"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.91%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = "'Polkadot"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'4.987"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.21%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'1.801.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.28%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.6687"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.82%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'83.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.43%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'6.114"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.35%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'29.165.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.02%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000008266"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.41%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'225.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.45%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'12.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.55%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +0.05%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'7.132"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.28%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'1.0000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.02%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'160.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.70%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'8.639"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.26%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.1394"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.53%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'17.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.65%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'1.509"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.29%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'4.114"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.68%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'4.030"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.66%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.183"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.13%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.05311"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.01%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.869"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.89%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.7533"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.56%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.36%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -0.35%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.296.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.20%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -1.75%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.720"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.04%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.9186"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.81%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.08623"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +15.36%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'5.955"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.11%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.007"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.77%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'102.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.33%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'RocketPoolETH"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1.970.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.80%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'BabyDogeCoin"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000122"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.18%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -0.67%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.769"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.02%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'63.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.61%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'Cronos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.05931"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.38%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'9.028"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.28%  "
$ws.Range("E51").Style = "Normal"

